# Update refreshed odds / snapshot timestamps in the Betfair Back/Lay sheet.
# Column index reference: F=6 Odd_H_Back, G=7 Odd_H_Lay, H=8 Odd_A_Back, I=9 Odd_A_Lay,
# J=10 Odd_D_Back, K=11 Odd_D_Lay, N=14 Odd_Under15_FT_Back, P=16 Odd_Under25_FT_Back,
# Q=17 Odd_Over25_FT_Back, R=18 Odd_Under35_FT_Back, S=19 Odd_Over35_FT_Back,
# T=20 Odd_BTTS_Yes_Back, U=21 Odd_BTTS_No_Back, X=24 Odd_CS_0x0_Lay, AB=28 Odd_CS_1x0_Lay,
# AC=29 Odd_CS_1x1_Lay, AH=34 Odd_CS_2x2_Lay, AI=35 Odd_CS_2x3_Lay, AJ=36 Odd_CS_3x0_Lay,
# AL=38 Odd_CS_3x2_Lay, AP=42 Odd_CS_0x0_Back, AQ=43 Odd_CS_0x1_Back, AR=44 Odd_CS_0x2_Back,
# AS=45 Odd_CS_0x3_Back, AT=46 Odd_CS_1x0_Back, AV=48 Odd_CS_1x2_Back, AW=49 Odd_CS_1x3_Back,
# AX=50 Odd_CS_2x0_Back, AZ=52 Odd_CS_2x2_Back, BA=53 Odd_CS_2x3_Back, BD=56 Odd_CS_3x2_Back,
# BE=57 Odd_CS_3x3_Back, BF=58 Odd_CS_Goleada_H_Back, BG=59 Odd_CS_Goleada_A_Back, BH=60 SnapshotTS

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 6).Value = 2.22
$ws.Cells.Item(2, 8).Value = 3.3
$ws.Cells.Item(2, 9).Value = 3.35
$ws.Cells.Item(2, 14).Value = 4.5
$ws.Cells.Item(2, 19).Value = 2.92
$ws.Cells.Item(2, 20).Value = 1.68
$ws.Cells.Item(2, 52).Value = 15
$ws.Cells.Item(2, 56).Value = 28
$ws.Cells.Item(2, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(3, 6).Value = 2.84
$ws.Cells.Item(3, 7).Value = 3.1
$ws.Cells.Item(3, 8).Value = 2.64
$ws.Cells.Item(3, 11).Value = 3.55
$ws.Cells.Item(3, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(4, 6).Value = 12.5
$ws.Cells.Item(4, 7).Value = 19.5
$ws.Cells.Item(4, 11).Value = 12
$ws.Cells.Item(4, 16).Value = 3.6
$ws.Cells.Item(4, 17).Value = 1.32
$ws.Cells.Item(4, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(5, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(6, 8).Value = 3.95
$ws.Cells.Item(6, 11).Value = 3.4
$ws.Cells.Item(6, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(7, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(8, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(9, 7).Value = 3.15
$ws.Cells.Item(9, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(10, 17).Value = 1.94
$ws.Cells.Item(10, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(11, 8).Value = 1.97
$ws.Cells.Item(11, 9).Value = 2.3
$ws.Cells.Item(11, 11).Value = 4.5
$ws.Cells.Item(11, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(12, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(13, 10).Value = 3.05
$ws.Cells.Item(13, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(14, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(15, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(16, 6).Value = 1.55
$ws.Cells.Item(16, 7).Value = 1.56
$ws.Cells.Item(16, 8).Value = 6
$ws.Cells.Item(16, 9).Value = 6.2
$ws.Cells.Item(16, 11).Value = 5.2
$ws.Cells.Item(16, 14).Value = 6.4
$ws.Cells.Item(16, 42).Value = 27
$ws.Cells.Item(16, 44).Value = 44
$ws.Cells.Item(16, 45).Value = 44
$ws.Cells.Item(16, 52).Value = 16.5
$ws.Cells.Item(16, 59).Value = 28
$ws.Cells.Item(16, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(17, 9).Value = 7
$ws.Cells.Item(17, 10).Value = 5.2
$ws.Cells.Item(17, 11).Value = 5.4
$ws.Cells.Item(17, 18).Value = 1.75
$ws.Cells.Item(17, 20).Value = 1.7
$ws.Cells.Item(17, 34).Value = 20
$ws.Cells.Item(17, 48).Value = 22
$ws.Cells.Item(17, 53).Value = 50
$ws.Cells.Item(17, 56).Value = 22
$ws.Cells.Item(17, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(18, 7).Value = 1.31
$ws.Cells.Item(18, 8).Value = 11
$ws.Cells.Item(18, 11).Value = 7.4
$ws.Cells.Item(18, 14).Value = 9
$ws.Cells.Item(18, 17).Value = 1.35
$ws.Cells.Item(18, 18).Value = 2.08
$ws.Cells.Item(18, 19).Value = 1.87
$ws.Cells.Item(18, 20).Value = 1.71
$ws.Cells.Item(18, 21).Value = 2.36
$ws.Cells.Item(18, 24).Value = 48
$ws.Cells.Item(18, 28).Value = 16
$ws.Cells.Item(18, 29).Value = 17.5
$ws.Cells.Item(18, 34).Value = 25
$ws.Cells.Item(18, 35).Value = 90
$ws.Cells.Item(18, 36).Value = 12
$ws.Cells.Item(18, 38).Value = 26
$ws.Cells.Item(18, 42).Value = 42
$ws.Cells.Item(18, 43).Value = 50
$ws.Cells.Item(18, 44).Value = 95
$ws.Cells.Item(18, 45).Value = 55
$ws.Cells.Item(18, 46).Value = 15
$ws.Cells.Item(18, 48).Value = 36
$ws.Cells.Item(18, 49).Value = 40
$ws.Cells.Item(18, 50).Value = 10.5
$ws.Cells.Item(18, 52).Value = 22
$ws.Cells.Item(18, 53).Value = 75
$ws.Cells.Item(18, 56).Value = 23
$ws.Cells.Item(18, 57).Value = 60
$ws.Cells.Item(18, 58).Value = 3.15
$ws.Cells.Item(18, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(19, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(20, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(21, 6).Value = 1.82
$ws.Cells.Item(21, 7).Value = 2.02
$ws.Cells.Item(21, 8).Value = 4.7
$ws.Cells.Item(21, 10).Value = 2.78
$ws.Cells.Item(21, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(22, 60).Value = "2026-02-23 06:21:43"

$ws.Cells.Item(23, 60).Value = "2026-02-23 06:21:43"
